# Updates crypto price/volume snapshot values per the latest scrape.
# For D-column values that look like plain numbers, assign with a leading
# apostrophe so Excel keeps them as text (matching the source inline-string
# cells) instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: '25.868.12' -> '25.883.05'; E2: '  +0.61%  ' -> '  +0.66%  '
$ws.Range("D2").Value = '25.883.05'
$ws.Range("E2").Value = '  +0.66%  '

# Row 3: D3: '1.629.67' -> '1.630.42'; E3: '  +0.06%  ' -> '  +0.10%  '
$ws.Range("D3").Value = '1.630.42'
$ws.Range("E3").Value = '  +0.10%  '

# Row 4: E4: '  -0.01%  ' -> '  -0.04%  '
$ws.Range("E4").Value = '  -0.04%  '

# Row 5: D5: '214.68' -> '214.76'; E5: '  +0.20%  ' -> '  +0.24%  '
$ws.Range("D5").Value = "'214.76"
$ws.Range("E5").Value = '  +0.24%  '

# Row 6: E6: '  +0.15%  ' -> '  +0.22%  '
$ws.Range("E6").Value = '  +0.22%  '

# Row 7: E7: '  -0.02%  ' -> '  -0.07%  '
$ws.Range("E7").Value = '  -0.07%  '

# Row 8: E8: '  +0.19%  ' -> '  +0.08%  '
$ws.Range("E8").Value = '  +0.08%  '

# Row 9: D9: '0.0631' -> '0.0632'; E9: '  +0.04%  ' -> '  +0.11%  '
$ws.Range("D9").Value = "'0.0632"
$ws.Range("E9").Value = '  +0.11%  '

# Row 10: E10: '  +1.05%  ' -> '  +0.88%  '
$ws.Range("E10").Value = '  +0.88%  '

# Row 11: D11: '0.0786' -> '0.0787'; E11: '  -0.45%  ' -> '  -0.63%  '
$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = '  -0.63%  '

# Row 12: D12: '1.855.00' -> '1.855.60'
$ws.Range("D12").Value = '1.855.60'

# Row 13: E13: '  -0.39%  ' -> '  -0.45%  '
$ws.Range("E13").Value = '  -0.45%  '

# Row 14: D14: '1.640.87' -> '1.622.58'; E14: '  -0.41%  ' -> '  -0.47%  '
$ws.Range("D14").Value = '1.622.58'
$ws.Range("E14").Value = '  -0.47%  '

# Row 15: E15: '  -1.99%  ' -> '  -1.94%  '
$ws.Range("E15").Value = '  -1.94%  '

# Row 16: D16: '0.0₃0757' -> '0.0₃0758'; E16: '  -0.28%  ' -> '  -0.25%  '
$ws.Range("D16").Value = '0.0₃0758'
$ws.Range("E16").Value = '  -0.25%  '

# Row 17: D17: '62.74' -> '62.77'; E17: '  -0.25%  ' -> '  -0.26%  '
$ws.Range("D17").Value = "'62.77"
$ws.Range("E17").Value = '  -0.26%  '

# Row 18: D18: '25.860.20' -> '25.864.72'
$ws.Range("D18").Value = '25.864.72'

# Row 19: E19: '  -0.07%  ' -> '  -0.03%  '
$ws.Range("E19").Value = '  -0.03%  '

# Row 20: E20: '  -1.28%  ' -> '  -1.19%  '
$ws.Range("E20").Value = '  -1.19%  '

# Row 21: D21: '192.68' -> '192.69'; E21: '  +0.70%  ' -> '  +0.60%  '
$ws.Range("D21").Value = "'192.69"
$ws.Range("E21").Value = '  +0.60%  '

# Row 22: E22: '  +0.80%  ' -> '  +0.62%  '
$ws.Range("E22").Value = '  +0.62%  '

# Row 23: D23: '6.25' -> '6.26'; E23: '  +0.09%  ' -> '  +0.27%  '
$ws.Range("D23").Value = "'6.26"
$ws.Range("E23").Value = '  +0.27%  '

# Row 24: E24: '  -2.26%  ' -> '  -2.27%  '
$ws.Range("E24").Value = '  -2.27%  '

# Row 25: B25: 'Monero' -> 'BinanceUSD'; C25: 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' -> 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D25: '142.99' -> '1.00'; E25: '  +0.54%  ' -> '  -0.15%  '
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = '  -0.15%  '

# Row 26: B26: 'BinanceUSD' -> 'Monero'; C26: 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' -> 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D26: '1.00' -> '142.58'; E26: '  -0.11%  ' -> '  +0.25%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'142.58"
$ws.Range("E26").Value = '  +0.25%  '

# Row 27: E27: '  +1.45%  ' -> '  +1.60%  '
$ws.Range("E27").Value = '  +1.60%  '

# Row 28: E28: '  +0.05%  ' -> '  +0.07%  '
$ws.Range("E28").Value = '  +0.07%  '

# Row 29: D29: '15.45' -> '15.44'; E29: '  +0.04%  ' -> '  -0.04%  '
$ws.Range("D29").Value = "'15.44"
$ws.Range("E29").Value = '  -0.04%  '

# Row 31: E31: '  +2.03%  ' -> '  +2.06%  '
$ws.Range("E31").Value = '  +2.06%  '

# Row 32: D32: '3.31' -> '3.30'; E32: '  -0.45%  ' -> '  -0.65%  '
$ws.Range("D32").Value = "'3.30"
$ws.Range("E32").Value = '  -0.65%  '

# Row 33: E33: '  -0.17%  ' -> '  -0.20%  '
$ws.Range("E33").Value = '  -0.20%  '

# Row 34: E34: '  +0.75%  ' -> '  +0.50%  '
$ws.Range("E34").Value = '  +0.50%  '

# Row 35: E35: '  +1.61%  ' -> '  +0.87%  '
$ws.Range("E35").Value = '  +0.87%  '

# Row 36: E36: '  -0.22%  ' -> '  -0.16%  '
$ws.Range("E36").Value = '  -0.16%  '

# Row 37: D37: '1.136.60' -> '1.136.04'; E37: '  +0.23%  ' -> '  +0.29%  '
$ws.Range("D37").Value = '1.136.04'
$ws.Range("E37").Value = '  +0.29%  '

# Row 38: E38: '  +1.73%  ' -> '  +1.44%  '
$ws.Range("E38").Value = '  +1.44%  '

# Row 39: D39: '2.47' -> '2.46'; E39: '  -2.02%  ' -> '  -1.95%  '
$ws.Range("D39").Value = "'2.46"
$ws.Range("E39").Value = '  -1.95%  '

# Row 40: E40: '  +0.78%  ' -> '  +0.55%  '
$ws.Range("E40").Value = '  +0.55%  '

# Row 41: E41: '  -0.15%  ' -> '  -0.11%  '
$ws.Range("E41").Value = '  -0.11%  '

# Row 42: E42: '  -0.75%  ' -> '  -0.93%  '
$ws.Range("E42").Value = '  -0.93%  '

# Row 43: D43: '98.99' -> '99.17'; E43: '  -1.54%  ' -> '  -1.40%  '
$ws.Range("D43").Value = "'99.17"
$ws.Range("E43").Value = '  -1.40%  '

# Row 44: D44: '0.800' -> '0.802'; E44: '  +0.20%  ' -> '  +0.29%  '
$ws.Range("D44").Value = "'0.802"
$ws.Range("E44").Value = '  +0.29%  '

# Row 45: D45: '1.765.13' -> '1.765.64'; E45: '  +0.12%  ' -> '  +0.11%  '
$ws.Range("D45").Value = '1.765.64'
$ws.Range("E45").Value = '  +0.11%  '

# Row 46: D46: '0.0₆0111' -> '0.0₆0112'; E46: '  +0.08%  ' -> '  +0.31%  '
$ws.Range("D46").Value = '0.0₆0112'
$ws.Range("E46").Value = '  +0.31%  '

# Row 47: D47: '56.14' -> '56.11'; E47: '  +1.61%  ' -> '  +1.69%  '
$ws.Range("D47").Value = "'56.11"
$ws.Range("E47").Value = '  +1.69%  '

# Row 48: D48: '0.0530' -> '0.0531'; E48: '  +4.43%  ' -> '  +4.66%  '
$ws.Range("D48").Value = "'0.0531"
$ws.Range("E48").Value = '  +4.66%  '

# Row 49: E49: '  +1.74%  ' -> '  +1.90%  '
$ws.Range("E49").Value = '  +1.90%  '

# Row 51: D51: '7.62' -> '7.63'; E51: '  +2.94%  ' -> '  +2.78%  '
$ws.Range("D51").Value = "'7.63"
$ws.Range("E51").Value = '  +2.78%  '

